# Restore the default "Office Theme" colour scheme on the presentation's
# theme (the deck had the "Integral" theme colours applied; the commit
# reverts the live theme part back to the stock Office Theme palette).
#
# Office Theme colour scheme (12 slots, in the standard
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order):
#   dk1      000000
#   lt1      FFFFFF
#   dk2      44546A
#   lt2      E7E6E6
#   accent1  5B9BD5
#   accent2  ED7D31
#   accent3  A5A5A5
#   accent4  FFC000
#   accent5  4472C4
#   accent6  70AD47
#   hlink    0563C1
#   folHlink 954F72

function ConvertTo-ComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$theme = $sm.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-ComRgb $officeThemeColors[$i - 1]
}
